# Update "want to go" counts (column F) for a set of events (identified by unique
# Bilibili "id=" fragment in the Link column H) and remove one finished event
# row ("杭州·懒喵N²..." on 2024-08-17, id=89522) from the sheets that list it
# ("演出" and "全部类型"), then renumber the index column (A) for the sheets
# whose row count changed.
#
# NOTE: this PowerShell runtime does not support named parameter binding
# ( -foo bar ) for user-defined functions, so all helper functions below use
# purely positional parameters.

function Update-WantToGoCounts($ws, $updates) {
    foreach ($key in $updates.Keys) {
        $found = $ws.Columns.Item(8).Find($key)
        if ($found -ne $null) {
            $r = $found.Row
            $ws.Cells.Item($r, 6).Value = $updates[$key]
        }
    }
}

function Remove-EventRow($ws, $key) {
    $found = $ws.Columns.Item(8).Find($key)
    if ($found -ne $null) {
        $r = $found.Row
        $ws.Rows.Item($r).Delete() | Out-Null
    }
}

function Renumber-IdColumn($ws) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}

$wb = $excel.ActiveWorkbook

# Updated "want to go" counts shared across sheets (keyed by unique Link id fragment)
$wantToGoUpdates = @{
    "id=89879" = 384
    "id=84912" = 1076
    "id=90773" = 25
    "id=87293" = 1042
    "id=87230" = 518
    "id=90025" = 401
    "id=90372" = 295
    "id=88899" = 333
    "id=90433" = 318
    "id=88498" = 417
    "id=83822" = 5385
    "id=89250" = 1507
    "id=89550" = 341
    "id=90057" = 4474
    "id=89966" = 1429
    "id=88452" = 628
    "id=88429" = 3769
    "id=86604" = 9347
    "id=89180" = 2095
}

$removedEventKey = "id=89522"

# 展览 (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
Update-WantToGoCounts $wsExhibition $wantToGoUpdates

# 演出 (Performance) - remove the finished "懒喵N²" event and shift rows up
$wsPerformance = $wb.Worksheets.Item("演出")
Update-WantToGoCounts $wsPerformance $wantToGoUpdates
Remove-EventRow $wsPerformance $removedEventKey
Renumber-IdColumn $wsPerformance

# 本地生活 (Local life)
$wsLocal = $wb.Worksheets.Item("本地生活")
Update-WantToGoCounts $wsLocal $wantToGoUpdates

# 全部类型 (All types) - union sheet: also remove the finished event row
$wsAll = $wb.Worksheets.Item("全部类型")
Update-WantToGoCounts $wsAll $wantToGoUpdates
Remove-EventRow $wsAll $removedEventKey
Renumber-IdColumn $wsAll
